$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Match the bold/border/center-top-alignment style used by the other header
# cells (e.g. AC1) by cloning its format onto the new header cells first...
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

# ...then set the new header labels for the new columns AD, AE, AF (row 1)
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Fill in team record (Wins/Losses/Ties) for every data row (2..57)
for ($row = 2; $row -le 57; $row++) {
    $ws.Cells.Item($row, 30).Value = 59
    $ws.Cells.Item($row, 31).Value = 103
    $ws.Cells.Item($row, 32).Value = 0
}
